# Update the canonical terminology URLs and the generation date
# in the DiplomeEtat StructureDefinition workbook.

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsElements = $wb.Worksheets.Item("Elements")

# Metadata sheet: refresh the "Date" row value (B8)
$wsMetadata.Range("B8").Value = "2025-07-25T07:22:51+00:00"

# Elements sheet: update canonical terminology URLs (column Z - "Binding Value Set")
$wsElements.Range("Z3").Value = "https://mos.esante.gouv.fr/NOS/TRE_R14-TypeDiplome/FHIR/TRE-R14-TypeDiplome?vs"
$wsElements.Range("Z4").Value = "https://mos.esante.gouv.fr/NOS/TRE_R16-LieuFormation/FHIR/TRE-R16-LieuFormation?vs"
$wsElements.Range("Z7").Value = "https://mos.esante.gouv.fr/NOS/TRE_R48-DiplomeEtatFrancais/FHIR/TRE-R48-DiplomeEtatFrancais?vs"

# Column Z widened (best-fit) to accommodate the new, longer URLs
$wsElements.Columns("Z").ColumnWidth = 80
